$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 39.249740724287321
$ws.Range("C2").Value = 39.875510433793252
$ws.Range("D2").Value = 46.228756715520007
$ws.Range("E2").Value = 32.128865179222252

$ws.Range("B3").Value = 68.194223218324865
$ws.Range("C3").Value = 43.396973129585184
$ws.Range("D3").Value = 43.044488056068573
$ws.Range("E3").Value = 16.955872732733589

$ws.Range("B1:E3").Select()
